$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 293.18182
$ws.Cells.Item(12, 9).Value = 312.44446
$ws.Cells.Item(12, 10).Value = 206.5
$ws.Cells.Item(12, 11).Value = 312.44446
$ws.Cells.Item(12, 12).Value = 206.5
$ws.Cells.Item(12, 13).Value = -142.44446
$ws.Cells.Item(12, 14).Value = -546.5
$ws.Cells.Item(98, 8).Value = 3755.0667
$ws.Cells.Item(98, 9).Value = 3539
$ws.Cells.Item(98, 11).Value = 3539
$ws.Cells.Item(98, 13).Value = -2041
$ws.Cells.Item(107, 8).Value = 2363.125
$ws.Cells.Item(107, 9).Value = 1083.25
$ws.Cells.Item(107, 10).Value = 6202.75
$ws.Cells.Item(107, 11).Value = 1083.25
$ws.Cells.Item(107, 12).Value = 6202.75
$ws.Cells.Item(107, 13).Value = 836.75
$ws.Cells.Item(107, 14).Value = -10042.75
$ws.Cells.Item(122, 8).Value = 3755.0667
$ws.Cells.Item(122, 9).Value = 3539
$ws.Cells.Item(122, 11).Value = 10617
$ws.Cells.Item(122, 13).Value = -8167
$ws.Cells.Item(135, 8).Value = 2250
$ws.Cells.Item(135, 9).Value = 1500
$ws.Cells.Item(135, 10).Value = 3000
$ws.Cells.Item(135, 11).Value = 13500
$ws.Cells.Item(135, 12).Value = 27000
$ws.Cells.Item(135, 13).Value = -10965
$ws.Cells.Item(135, 14).Value = -32070
$ws.Cells.Item(138, 8).Value = 1679.7273
$ws.Cells.Item(138, 9).Value = 909.625
$ws.Cells.Item(138, 11).Value = 2728.875
$ws.Cells.Item(138, 13).Value = 2411.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(62, 8).Value = 49999.5
$ws.Cells.Item(62, 10).Value = 49999.5
$ws.Cells.Item(62, 12).Value = 49999.5
$ws.Cells.Item(62, 14).Value = -51247.5
$ws.Cells.Item(65, 8).Value = 49999.5
$ws.Cells.Item(65, 10).Value = 49999.5
$ws.Cells.Item(65, 12).Value = 149998.5
$ws.Cells.Item(65, 14).Value = -156238.5
$ws.Cells.Item(74, 8).Value = 1476.2
$ws.Cells.Item(74, 9).Value = 1185.7273
$ws.Cells.Item(74, 10).Value = 2275
$ws.Cells.Item(74, 11).Value = 1185.7273
$ws.Cells.Item(74, 12).Value = 2275
$ws.Cells.Item(74, 13).Value = -311.7273
$ws.Cells.Item(74, 14).Value = -4023
$ws.Cells.Item(77, 8).Value = 1476.2
$ws.Cells.Item(77, 9).Value = 1185.7273
$ws.Cells.Item(77, 10).Value = 2275
$ws.Cells.Item(77, 11).Value = 5928.636500000001
$ws.Cells.Item(77, 12).Value = 11375
$ws.Cells.Item(77, 13).Value = -1560.636500000001
$ws.Cells.Item(77, 14).Value = -20111
$ws.Cells.Item(110, 8).Value = 3603.6428
$ws.Cells.Item(110, 9).Value = 1490.2
$ws.Cells.Item(110, 11).Value = 1490.2
$ws.Cells.Item(110, 13).Value = 554.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3103.6072
$ws.Cells.Item(86, 9).Value = 3118.1
$ws.Cells.Item(86, 10).Value = 3067.375
$ws.Cells.Item(86, 11).Value = 3118.1
$ws.Cells.Item(86, 12).Value = 3067.375
$ws.Cells.Item(86, 13).Value = -1995.1
$ws.Cells.Item(86, 14).Value = -5313.375
$ws.Cells.Item(89, 8).Value = 3103.6072
$ws.Cells.Item(89, 9).Value = 3118.1
$ws.Cells.Item(89, 10).Value = 3067.375
$ws.Cells.Item(89, 11).Value = 15590.5
$ws.Cells.Item(89, 12).Value = 15336.875
$ws.Cells.Item(89, 13).Value = -9974.5
$ws.Cells.Item(89, 14).Value = -26568.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3712.6428
$ws.Cells.Item(31, 9).Value = 2715.2
$ws.Cells.Item(31, 11).Value = 2715.2
$ws.Cells.Item(31, 13).Value = -2420.2
$ws.Cells.Item(34, 8).Value = 3712.6428
$ws.Cells.Item(34, 9).Value = 2715.2
$ws.Cells.Item(34, 11).Value = 2715.2
$ws.Cells.Item(34, 13).Value = -2513.2
$ws.Cells.Item(107, 8).Value = 335.73334
$ws.Cells.Item(107, 10).Value = 447.83334
$ws.Cells.Item(107, 12).Value = 447.83334
$ws.Cells.Item(107, 14).Value = -4287.83334
$ws.Cells.Item(132, 8).Value = 2514.1428
$ws.Cells.Item(132, 9).Value = 2514.1428
$ws.Cells.Item(132, 11).Value = 7542.428400000001
$ws.Cells.Item(132, 13).Value = -5012.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(110, 8).Value = 2900
$ws.Cells.Item(110, 9).Value = 2900
$ws.Cells.Item(110, 11).Value = 8700
$ws.Cells.Item(110, 13).Value = -4610
$ws.Cells.Item(113, 8).Value = 1419.7
$ws.Cells.Item(113, 10).Value = 1521.8889
$ws.Cells.Item(113, 12).Value = 4565.6667
$ws.Cells.Item(113, 14).Value = -8905.6667
$ws.Cells.Item(129, 8).Value = 2209.4614
$ws.Cells.Item(129, 9).Value = 904.5
$ws.Cells.Item(129, 10).Value = 3328
$ws.Cells.Item(129, 11).Value = 2713.5
$ws.Cells.Item(129, 12).Value = 9984
$ws.Cells.Item(129, 13).Value = 2286.5
$ws.Cells.Item(129, 14).Value = -19984
$ws.Cells.Item(138, 8).Value = 1174.4
$ws.Cells.Item(138, 9).Value = 624
$ws.Cells.Item(138, 10).Value = 2000
$ws.Cells.Item(138, 11).Value = 1872
$ws.Cells.Item(138, 12).Value = 6000
$ws.Cells.Item(138, 13).Value = 3268
$ws.Cells.Item(138, 14).Value = -16280

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 13427
$ws.Cells.Item(97, 8).Value = 851.05884
$ws.Cells.Item(97, 9).Value = 851.05884
$ws.Cells.Item(97, 11).Value = 851.05884
$ws.Cells.Item(97, 13).Value = -355.05884
$ws.Cells.Item(102, 8).Value = 1611.4166
$ws.Cells.Item(102, 9).Value = 1611.4166
$ws.Cells.Item(102, 11).Value = 1611.4166
$ws.Cells.Item(102, 13).Value = 10.58339999999998
$ws.Cells.Item(126, 8).Value = 3666.6667
$ws.Cells.Item(126, 9).Value = 3666.6667
$ws.Cells.Item(126, 11).Value = 11000.0001
$ws.Cells.Item(126, 13).Value = -8530.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1408.5
$ws.Cells.Item(22, 9).Value = 1116
$ws.Cells.Item(22, 10).Value = 1847.25
$ws.Cells.Item(22, 11).Value = 1116
$ws.Cells.Item(22, 12).Value = 1847.25
$ws.Cells.Item(22, 13).Value = -821
$ws.Cells.Item(22, 14).Value = -2437.25
$ws.Cells.Item(27, 8).Value = 1408.5
$ws.Cells.Item(27, 9).Value = 1116
$ws.Cells.Item(27, 10).Value = 1847.25
$ws.Cells.Item(27, 11).Value = 1116
$ws.Cells.Item(27, 12).Value = 1847.25
$ws.Cells.Item(27, 13).Value = -1009
$ws.Cells.Item(27, 14).Value = -2061.25
$ws.Cells.Item(82, 8).Value = 1604.25
$ws.Cells.Item(82, 9).Value = 2025.3334
$ws.Cells.Item(82, 10).Value = 1183.1666
$ws.Cells.Item(82, 11).Value = 2025.3334
$ws.Cells.Item(82, 12).Value = 1183.1666
$ws.Cells.Item(82, 13).Value = -1664.3334
$ws.Cells.Item(82, 14).Value = -1905.1666
$ws.Cells.Item(85, 8).Value = 1604.25
$ws.Cells.Item(85, 9).Value = 2025.3334
$ws.Cells.Item(85, 10).Value = 1183.1666
$ws.Cells.Item(85, 11).Value = 2025.3334
$ws.Cells.Item(85, 12).Value = 1183.1666
$ws.Cells.Item(85, 13).Value = -777.3334
$ws.Cells.Item(85, 14).Value = -3679.1666
$ws.Cells.Item(132, 8).Value = 19537.092
$ws.Cells.Item(132, 9).Value = 20858.285
$ws.Cells.Item(132, 10).Value = 17225
$ws.Cells.Item(132, 11).Value = 62574.855
$ws.Cells.Item(132, 12).Value = 51675
$ws.Cells.Item(132, 13).Value = -60044.855
$ws.Cells.Item(132, 14).Value = -56735
$ws.Cells.Item(136, 8).Value = 3512.5833
$ws.Cells.Item(136, 9).Value = 3292.3076
$ws.Cells.Item(136, 10).Value = 3772.9092
$ws.Cells.Item(136, 11).Value = 9876.9228
$ws.Cells.Item(136, 12).Value = 11318.7276
$ws.Cells.Item(136, 13).Value = -7326.9228
$ws.Cells.Item(136, 14).Value = -16418.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1214.7222
$ws.Cells.Item(113, 9).Value = 825.6
$ws.Cells.Item(113, 10).Value = 1701.125
$ws.Cells.Item(113, 11).Value = 2476.8
$ws.Cells.Item(113, 12).Value = 5103.375
$ws.Cells.Item(113, 13).Value = -306.8000000000002
$ws.Cells.Item(113, 14).Value = -9443.375
$ws.Cells.Item(136, 8).Value = 2407.8
$ws.Cells.Item(136, 9).Value = 2407.8
$ws.Cells.Item(136, 11).Value = 7223.400000000001
$ws.Cells.Item(136, 13).Value = -4673.400000000001
